# Updated cryptos list values (Price / Volume(1h) columns), applied as plain text
# to preserve the exact display strings (e.g. thousand-separator dots, padded %).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '71.105.46'
$c = $ws.Range("E2")
$c.NumberFormat = '@'
$c.Value = '  +3.14%  '
$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '3.611.80'
$c = $ws.Range("E3")
$c.NumberFormat = '@'
$c.Value = '  +2.73%  '
$c = $ws.Range("E4")
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '601.71'
$c = $ws.Range("E5")
$c.NumberFormat = '@'
$c.Value = '  +2.64%  '
$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '175.02'
$c = $ws.Range("E6")
$c.NumberFormat = '@'
$c.Value = '  +2.23%  '
$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '3.603.62'
$c = $ws.Range("E7")
$c.NumberFormat = '@'
$c.Value = '  +2.68%  '
$c = $ws.Range("D8")
$c.NumberFormat = '@'
$c.Value = '0.619'
$c = $ws.Range("E8")
$c.NumberFormat = '@'
$c.Value = '  +1.42%  '
$c = $ws.Range("E9")
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '
$c = $ws.Range("E10")
$c.NumberFormat = '@'
$c.Value = '  +7.04%  '
$c = $ws.Range("E11")
$c.NumberFormat = '@'
$c.Value = '  +7.51%  '
$c = $ws.Range("E12")
$c.NumberFormat = '@'
$c.Value = '  +2.36%  '
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '47.35'
$c = $ws.Range("E13")
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '
$c = $ws.Range("E14")
$c.NumberFormat = '@'
$c.Value = '  +2.00%  '
$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '4.192.24'
$c = $ws.Range("E15")
$c.NumberFormat = '@'
$c.Value = '  +2.93%  '
$c = $ws.Range("E16")
$c.NumberFormat = '@'
$c.Value = '  +0.47%  '
$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '620.09'
$c = $ws.Range("E17")
$c.NumberFormat = '@'
$c.Value = '  -1.18%  '
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '3.620.34'
$c = $ws.Range("E18")
$c.NumberFormat = '@'
$c.Value = '  +2.62%  '
$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '71.238.80'
$c = $ws.Range("E19")
$c.NumberFormat = '@'
$c.Value = '  +3.23%  '
$c = $ws.Range("E20")
$c.NumberFormat = '@'
$c.Value = '  -1.35%  '
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '17.60'
$c = $ws.Range("E21")
$c.NumberFormat = '@'
$c.Value = '  +1.14%  '
$c = $ws.Range("E22")
$c.NumberFormat = '@'
$c.Value = '  +0.99%  '
$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '9.36'
$c = $ws.Range("E23")
$c.NumberFormat = '@'
$c.Value = '  -16.02%  '
$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '16.17'
$c = $ws.Range("E24")
$c.NumberFormat = '@'
$c.Value = '  +1.45%  '
$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '98.15'
$c = $ws.Range("E25")
$c.NumberFormat = '@'
$c.Value = '  +1.24%  '
$c = $ws.Range("E26")
$c.NumberFormat = '@'
$c.Value = '  -0.45%  '
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '2.68'
$c = $ws.Range("E27")
$c.NumberFormat = '@'
$c.Value = '  +2.28%  '
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '0.999'
$c = $ws.Range("E28")
$c.NumberFormat = '@'
$c.Value = '  -0.06%  '
$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '34.34'
$c = $ws.Range("E29")
$c.NumberFormat = '@'
$c.Value = '  +4.81%  '
$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '9.35'
$c = $ws.Range("E30")
$c.NumberFormat = '@'
$c.Value = '  +0.97%  '
$c = $ws.Range("E31")
$c.NumberFormat = '@'
$c.Value = '  +0.34%  '
$c = $ws.Range("E32")
$c.NumberFormat = '@'
$c.Value = '  -0.99%  '
$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '7.30'
$c = $ws.Range("E33")
$c.NumberFormat = '@'
$c.Value = '  +5.29%  '
$c = $ws.Range("E34")
$c.NumberFormat = '@'
$c.Value = '  -0.58%  '
$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '629.79'
$c = $ws.Range("E35")
$c.NumberFormat = '@'
$c.Value = '  -1.51%  '
$c = $ws.Range("D36")
$c.NumberFormat = '@'
$c.Value = '3.77'
$c = $ws.Range("E36")
$c.NumberFormat = '@'
$c.Value = '  +8.44%  '
$c = $ws.Range("E37")
$c.NumberFormat = '@'
$c.Value = '  +0.03%  '
$c = $ws.Range("E38")
$c.NumberFormat = '@'
$c.Value = '  +1.73%  '
$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '0.0484'
$c = $ws.Range("E39")
$c.NumberFormat = '@'
$c.Value = '  +6.61%  '
$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '57.62'
$c = $ws.Range("E40")
$c.NumberFormat = '@'
$c.Value = '  +0.68%  '
$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '0.144'
$c = $ws.Range("E42")
$c.NumberFormat = '@'
$c.Value = '  +6.08%  '
$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '3.416.80'
$c = $ws.Range("E43")
$c.NumberFormat = '@'
$c.Value = '  +0.86%  '
$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '0.328'
$c = $ws.Range("E44")
$c.NumberFormat = '@'
$c.Value = '  +0.29%  '
$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '0.0₃0723'
$c = $ws.Range("E45")
$c.NumberFormat = '@'
$c.Value = '  +3.82%  '
$c = $ws.Range("E46")
$c.NumberFormat = '@'
$c.Value = '  +9.89%  '
$c = $ws.Range("B47")
$c.NumberFormat = '@'
$c.Value = 'InjectiveProtocol'
$c = $ws.Range("C47")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '33.33'
$c = $ws.Range("E47")
$c.NumberFormat = '@'
$c.Value = '  +1.64%  '
$c = $ws.Range("B48")
$c.NumberFormat = '@'
$c.Value = 'Fetch.AI'
$c = $ws.Range("C48")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '2.71'
$c = $ws.Range("E48")
$c.NumberFormat = '@'
$c.Value = '  +6.53%  '
$c = $ws.Range("E49")
$c.NumberFormat = '@'
$c.Value = '  +1.18%  '
$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '132.90'
$c = $ws.Range("E50")
$c.NumberFormat = '@'
$c.Value = '  +0.32%  '
$c = $ws.Range("E51")
$c.NumberFormat = '@'
$c.Value = '  -0.05%  '
